$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 ("Sponsor Work"): add D13 = 2, and I13 (Daily Total) updates from 2 to 4
$ws.Range("D13").Value = 2
$ws.Range("I13").Value = 4

# Row 14 ("Daily Total" row): D14 becomes 2 (was empty), and I14 (Weekly Total) updates from 2 to 4
$ws.Range("D14").Value = 2
$ws.Range("I14").Value = 4

# Update the active selection to K14
$ws.Range("K14").Select()
